$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("L2").Value = 1676
$ws.Range("K3").Value = 8182
$ws.Range("L3").Value = 1693
$ws.Range("I4").Value = 1832
$ws.Range("J4").Value = 1862
$ws.Range("K4").Value = 1756
$ws.Range("L4").Value = 478
$ws.Range("L6").Value = 1608
$ws.Range("I7").Value = 26299
$ws.Range("J7").Value = 29335
$ws.Range("L7").Value = 5558

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("L2").Value = 96
$ws.Range("L3").Value = 107
$ws.Range("L6").Value = 95
$ws.Range("L7").Value = 335

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("L2").Value = 38
$ws.Range("L7").Value = 129

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("L2").Value = 50
$ws.Range("L7").Value = 235

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("L2").Value = 61
$ws.Range("L3").Value = 55
$ws.Range("L7").Value = 195

$ws = $wb.Worksheets.Item("New City")
$ws.Range("L2").Value = 41
$ws.Range("L7").Value = 109

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("L2").Value = 20
$ws.Range("L7").Value = 85

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("L3").Value = 3
$ws.Range("L4").Value = 20
$ws.Range("L7").Value = 185
$ws.Range("L8").Value = 335
$ws.Range("L11").Value = 100
$ws.Range("L19").Value = 160
$ws.Range("L20").Value = 150
$ws.Range("L22").Value = 19
$ws.Range("L27").Value = 59
$ws.Range("L33").Value = 235
$ws.Range("L36").Value = 85
$ws.Range("L37").Value = 195
$ws.Range("L42").Value = 173
$ws.Range("L44").Value = 40
$ws.Range("L45").Value = 8
$ws.Range("L48").Value = 82
$ws.Range("L51").Value = 68
$ws.Range("L52").Value = 115
$ws.Range("L54").Value = 124
$ws.Range("L60").Value = 32
$ws.Range("I63").Value = 254
$ws.Range("J63").Value = 213
$ws.Range("K63").Value = 85
$ws.Range("L63").Value = 21
$ws.Range("L64").Value = 39
$ws.Range("L65").Value = 109
$ws.Range("L75").Value = 24
$ws.Range("L76").Value = 61
$ws.Range("L77").Value = 35
$ws.Range("L79").Value = 152
$ws.Range("L83").Value = 129
$ws.Range("L84").Value = 58
$ws.Range("L85").Value = 290
$ws.Range("L87").Value = 16
$ws.Range("L89").Value = 67
$ws.Range("K90").Value = 262
$ws.Range("L90").Value = 55
$ws.Range("L91").Value = 70
$ws.Range("L93").Value = 29
$ws.Range("L94").Value = 71
$ws.Range("L99").Value = 85
$ws.Range("I101").Value = 26299
$ws.Range("J101").Value = 29335
$ws.Range("L101").Value = 5558

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("L3").Value = 22
$ws.Range("L7").Value = 58

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("L2").Value = 31
$ws.Range("L7").Value = 124

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("L3").Value = 18
$ws.Range("L4").Value = 22
$ws.Range("L7").Value = 82

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("L3").Value = 51
$ws.Range("L6").Value = 52
$ws.Range("L7").Value = 160

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("L2").Value = 14
$ws.Range("L7").Value = 40

$ws = $wb.Worksheets.Item("River North")
$ws.Range("L6").Value = 30
$ws.Range("L7").Value = 61

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("L3").Value = 44
$ws.Range("L7").Value = 173

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("L3").Value = 23
$ws.Range("L7").Value = 70

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("L2").Value = 55
$ws.Range("L7").Value = 152

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("L3").Value = 8
$ws.Range("L6").Value = 11
$ws.Range("L7").Value = 39

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("L3").Value = 41
$ws.Range("L6").Value = 46
$ws.Range("L7").Value = 150

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("L4").Value = 6
$ws.Range("L7").Value = 85

$ws = $wb.Worksheets.Item("West Lawn")
$ws.Range("L6").Value = 11
$ws.Range("L7").Value = 29

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("L3").Value = 59
$ws.Range("L4").Value = 16
$ws.Range("L6").Value = 51
$ws.Range("L7").Value = 185

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("L3").Value = 17
$ws.Range("L6").Value = 23
$ws.Range("L7").Value = 71

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("L2").Value = 34
$ws.Range("L3").Value = 29
$ws.Range("L7").Value = 100

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("L3").Value = 17
$ws.Range("L7").Value = 67

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("L2").Value = 15
$ws.Range("L7").Value = 59

$ws = $wb.Worksheets.Item("Pullman")
$ws.Range("L2").Value = 13
$ws.Range("L7").Value = 24

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("L2").Value = 22
$ws.Range("K3").Value = 70
$ws.Range("L4").Value = 4
$ws.Range("K7").Value = 262
$ws.Range("L7").Value = 55

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("L4").Value = 9
$ws.Range("L7").Value = 68

$ws = $wb.Worksheets.Item("Morgan Park")
$ws.Range("L3").Value = 11
$ws.Range("L7").Value = 32

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("L2").Value = 85
$ws.Range("L4").Value = 23
$ws.Range("L6").Value = 56
$ws.Range("L7").Value = 290

$ws = $wb.Worksheets.Item("Clearing")
$ws.Range("L2").Value = 8
$ws.Range("L7").Value = 19

$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("L2").Value = 11
$ws.Range("L7").Value = 35

$ws = $wb.Worksheets.Item("Jackson Park")
$ws.Range("L4").Value = 1
$ws.Range("L7").Value = 8

$ws = $wb.Worksheets.Item("Andersonville")
$ws.Range("L2").Value = 1
$ws.Range("L7").Value = 3

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("L3").Value = 32
$ws.Range("L6").Value = 35
$ws.Range("L7").Value = 115

$ws = $wb.Worksheets.Item("Archer Heights")
$ws.Range("L3").Value = 7
$ws.Range("L7").Value = 20

$ws = $wb.Worksheets.Item("Ukrainian Village")
$ws.Range("L2").Value = 3
$ws.Range("L7").Value = 16
